$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-31 18:55:50"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H3").Value = "2016-08-31 18:55:45"
$wsZh.Range("K3").Value = "2016-08-31 18:56:09"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H3").Value = "2016-08-31 18:55:50"
$wsDe.Range("K3").Value = "2016-08-31 18:56:18"
